# Aragon hospital coronavirus occupancy data: add the day's new block of
# rows (date serial 43990 = 2020-06-08) for the 20 hospitals/centres,
# appended right after the existing last block (rows 1293:1312, which
# holds the previous day's data for date 43989 = 2020-06-07).
#
# The new block reuses the exact same per-row formatting/styles as the
# prior day's block (banding styles 2-6 alternate per hospital row), so
# the simplest and most faithful way to reproduce it is to copy that
# range down, then patch: (a) the date column to the new day, and
# (b) the handful of occupancy counts (camas_ocupadas_total /
# camas_uci_ocupadas) that actually differ from the previous day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the previous day's 20-row block (A1293:H1312) into the new
# 20-row block (A1313:H1332), carrying over values + styles/number
# formats exactly as Excel would when copy/pasting the block down.
$src = $ws.Range("A1293:H1312")
$dst = $ws.Range("A1313:H1332")
$src.Copy($dst)

# Update the date column for the newly appended block to the new day.
$ws.Range("A1313:A1332").Value = 43990

# Patch the cells whose values changed versus the previous day's block.
# row 1314 -> Hospital Clínico Universitario: camas_ocupadas_total 26 -> 28
$ws.Range("C1314").Value = 28
# row 1320 -> Hospital de Barbastro: camas_ocupadas_total 8 -> 10
$ws.Range("C1320").Value = 10
# row 1321 -> Hospital San Jorge: camas_ocupadas_total 4 -> 6
$ws.Range("C1321").Value = 6
# row 1328 -> Clínica Montpellier: camas_ocupadas_total (blank) -> 1, camas_uci_ocupadas (blank) -> 1
$ws.Range("C1328").Value = 1
$ws.Range("D1328").Value = 1
